# Weekly refresh of the "Poroto granado" sheet:
# a new week's record is added at row 48, and every existing record from
# row 48 down to row 148 shifts down by one row (row 149 being added to
# hold what used to be row 148's data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 48:148 down to 49:149, carrying their formatting with them.
$ws.Rows(48).Insert()

# Populate the now-empty row 48 with the new weekly record.
$ws.Cells.Item(48, 1).Value = 5
$ws.Cells.Item(48, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(48, 3).Value = "Maule"
$ws.Cells.Item(48, 4).Value = 44665
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value = 7
$ws.Cells.Item(48, 6).Value = 100112030
$ws.Cells.Item(48, 7).Value = "Poroto granado"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 200
$ws.Cells.Item(48, 11).Value = 20000
$ws.Cells.Item(48, 12).Value = 20000
$ws.Cells.Item(48, 13).Value = 20000
$ws.Cells.Item(48, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(48, 15).Value = "Región del Maule"
$ws.Cells.Item(48, 16).Value = 800
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
